$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old block of rows (10-24) and insert a fresh block of 14 rows (10-23)
# to cleanly reset content + row heights before repopulating.
$ws.Range("A10:C24").EntireRow.Delete()
$ws.Range("A10:C23").EntireRow.Insert()

# Repopulate cell values for rows 10-23 per the new layout
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C10").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Understand Operational Research as a science applied to Industrial Engineering. Provide knowledge of the typical problems encountered in Industrial Engineering. Analyze, model and solve problems through Operational Research."
$ws.Range("C11").Value = "Understand Operational Research as a science applied to Industrial Engineering. Provide knowledge of the typical problems encountered in Industrial Engineering. Analyze, model and solve problems through Operational Research."
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory."
$ws.Range("C14").Value = "Introduction to Operational Research, Linear Programming, Simplex Method, Introduction to Graphs and Network Optimization, Case Study in Linear Programming, Introduction to Queue Theory."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models"
$ws.Range("C16").Value = "1. Introduction to Operational Research 1.1. Concepts of Operational Research; 1.2. Modeling; 1.3. Structure of Mathematical Models; 1.4. Mathematical techniques in Operational Research; 1.2. Phases of a Study in Operational Research 2. Linear Programming 2.1. Definition 2.2. Formulation of Models 2.3. Graphic Resolution; 3. Simplex method 3.1. Development of the Simplex Method; 3.2. Simplex Method Procedure; 4. Introduction to Graphs and Network Optimization 4.1. Basic Concepts in Graph Theory 4.2. Maximum Flow Problems; 4.3. Minimum Path Problems 5. Case Studies in Linear Programming 5.1. Simple Transport Model 5.2. Model of Designation. 6. Introduction to Queuing Theory 6.1. Queuing Theory Concepts 6.2. Markovian Models"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("C18").Value = "5840917 - Fabrício Maciel Gomes"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"

# Set custom row heights per the new layout
$ws.Rows(10).RowHeight = 60
$ws.Rows(11).RowHeight = 60
$ws.Rows(13).RowHeight = 60
$ws.Rows(14).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(16).RowHeight = 120
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(23).RowHeight = 30
